$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (Name / Annual Income, mln RUR) below the
# existing table, which currently ends at row 141.
$ws.Range("A142").Value = "Архипов Станислав "
$ws.Range("B142").Value = 23

$ws.Range("A143").Value = "Струна Марина "
$ws.Range("B143").Value = 20

$ws.Range("A144").Value = "Седов Станислав"
$ws.Range("B144").Value = 17

# Scroll/select to show the newly added rows, matching the author's view state.
$ws.Range("A142:B144").Select()
$ws.Application.ActiveWindow.ScrollRow = 130
